$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new machine record (row 33) following the existing table pattern:
# id | name | mac_address | serial_num | ip_address | mspec_id | lang_code | is_active | cr_by | cr_dtimes | eff_dtimes
$ws.Range("A33").Value = 10032
$ws.Range("B33").Value = "Machine 32"
$ws.Range("C33").Value = "F4-30-B9-D4-CD-6F"
$ws.Range("D33").Value = "FB5962911665"
$ws.Range("E33").Value = "192.168.0.358"
$ws.Range("F33").Value = 1001
$ws.Range("G33").Value = "eng"
$ws.Range("H33").Value = $true
$ws.Range("I33").Value = "superadmin"
$ws.Range("J33").Value = "now()"
$ws.Range("K33").Value = "now()"

# Scroll the window so the newly added rows are in view and leave the
# active selection where the author left it.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C28").Select()
